$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 2728.9092
$ws.Range("I74").Value = 2651.8
$ws.Range("J74").Value = 3500
$ws.Range("K74").Value = 2651.8
$ws.Range("L74").Value = 3500
$ws.Range("M74").Value = -1715.8
$ws.Range("N74").Value = -5372
$ws.Range("H76").Value = 5100
$ws.Range("I76").Value = 4060.2
$ws.Range("K76").Value = 4060.2
$ws.Range("M76").Value = -3745.2
$ws.Range("H77").Value = 2728.9092
$ws.Range("I77").Value = 2651.8
$ws.Range("J77").Value = 3500
$ws.Range("K77").Value = 13259
$ws.Range("L77").Value = 17500
$ws.Range("M77").Value = -8579
$ws.Range("N77").Value = -26860
$ws.Range("H79").Value = 5100
$ws.Range("I79").Value = 4060.2
$ws.Range("K79").Value = 4060.2
$ws.Range("M79").Value = -2968.2
$ws.Range("H80").Value = 694
$ws.Range("J80").Value = 1167.25
$ws.Range("L80").Value = 3501.75
$ws.Range("N80").Value = -5497.75
$ws.Range("H83").Value = 694
$ws.Range("J83").Value = 1167.25
$ws.Range("L83").Value = 10505.25
$ws.Range("N83").Value = -20489.25
$ws.Range("H114").Value = 37722
$ws.Range("J114").Value = 37722
$ws.Range("L114").Value = 37722
$ws.Range("N114").Value = -46400
$ws.Range("H129").Value = 1099.0454
$ws.Range("I129").Value = 276.55554
$ws.Range("J129").Value = 1668.4615
$ws.Range("K129").Value = 829.66662
$ws.Range("L129").Value = 5005.3845
$ws.Range("M129").Value = 4170.33338
$ws.Range("N129").Value = -15005.3845
$ws.Range("H132").Value = 2749.8333
$ws.Range("I132").Value = 2328.8235
$ws.Range("J132").Value = 3772.2856
$ws.Range("K132").Value = 6986.470499999999
$ws.Range("L132").Value = 11316.8568
$ws.Range("M132").Value = -4456.470499999999
$ws.Range("N132").Value = -16376.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 17859538
$ws.Range("I61").Value = 23258084
$ws.Range("J61").Value = 2811.2307
$ws.Range("K61").Value = 23258084
$ws.Range("L61").Value = 2811.2307
$ws.Range("M61").Value = -23257872
$ws.Range("N61").Value = -3235.2307
$ws.Range("H88").Value = 5157.4
$ws.Range("I88").Value = 1800
$ws.Range("J88").Value = 5996.75
$ws.Range("K88").Value = 1800
$ws.Range("L88").Value = 5996.75
$ws.Range("M88").Value = -1394
$ws.Range("N88").Value = -6808.75
$ws.Range("H91").Value = 5157.4
$ws.Range("I91").Value = 1800
$ws.Range("J91").Value = 5996.75
$ws.Range("K91").Value = 1800
$ws.Range("L91").Value = 5996.75
$ws.Range("M91").Value = -396
$ws.Range("N91").Value = -8804.75
$ws.Range("H132").Value = 4777.7026
$ws.Range("I132").Value = 1563.2222
$ws.Range("J132").Value = 13456.8
$ws.Range("K132").Value = 4689.6666
$ws.Range("L132").Value = 40370.39999999999
$ws.Range("M132").Value = -2159.6666
$ws.Range("N132").Value = -45430.39999999999
$ws.Range("H136").Value = 17859538
$ws.Range("I136").Value = 23258084
$ws.Range("J136").Value = 2811.2307
$ws.Range("K136").Value = 69774252
$ws.Range("L136").Value = 8433.6921
$ws.Range("M136").Value = -69771702
$ws.Range("N136").Value = -13533.6921

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3006.1765
$ws.Range("I86").Value = 1831
$ws.Range("J86").Value = 4050.7778
$ws.Range("K86").Value = 1831
$ws.Range("L86").Value = 4050.7778
$ws.Range("M86").Value = -708
$ws.Range("N86").Value = -6296.7778
$ws.Range("H89").Value = 3006.1765
$ws.Range("I89").Value = 1831
$ws.Range("J89").Value = 4050.7778
$ws.Range("K89").Value = 9155
$ws.Range("L89").Value = 20253.889
$ws.Range("M89").Value = -3539
$ws.Range("N89").Value = -31485.889
$ws.Range("H99").Value = 2252.75
$ws.Range("J99").Value = 3505.5
$ws.Range("L99").Value = 3505.5
$ws.Range("N99").Value = -6501.5
$ws.Range("H134").Value = 1721
$ws.Range("I134").Value = 1681.7142
$ws.Range("J134").Value = 1776
$ws.Range("K134").Value = 5045.142599999999
$ws.Range("L134").Value = 5328
$ws.Range("M134").Value = -2510.142599999999
$ws.Range("N134").Value = -10398

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 184.91667
$ws.Range("J5").Value = 173.71428
$ws.Range("L5").Value = 173.71428
$ws.Range("N5").Value = -397.71428
$ws.Range("H14").Value = 1000
$ws.Range("I14").Value = 1000
$ws.Range("K14").Value = 1000
$ws.Range("M14").Value = -830
$ws.Range("H58").Value = 2021.7
$ws.Range("I58").Value = 2229
$ws.Range("K58").Value = 2229
$ws.Range("M58").Value = -2026
$ws.Range("H105").Value = 1083.1666
$ws.Range("I105").Value = 1083.1666
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 1083.1666
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = 663.8334
$ws.Range("N105").Value = ""
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").Value = ""
$ws.Range("H134").Value = 3168.889
$ws.Range("I134").Value = 3419.2727
$ws.Range("K134").Value = 10257.8181
$ws.Range("M134").Value = -7722.8181
$ws.Range("H136").Value = 2021.7
$ws.Range("I136").Value = 2229
$ws.Range("K136").Value = 6687
$ws.Range("M136").Value = -4137

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1084.9656
$ws.Range("I97").Value = 1087.5
$ws.Range("J97").Value = 1081.8462
$ws.Range("K97").Value = 1087.5
$ws.Range("L97").Value = 1081.8462
$ws.Range("M97").Value = -591.5
$ws.Range("N97").Value = -2073.8462
$ws.Range("H101").Value = 35657
$ws.Range("J101").Value = 35657
$ws.Range("L101").Value = 35657
$ws.Range("N101").Value = -42147

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1680.4
$ws.Range("I61").Value = 1126
$ws.Range("J61").Value = 2050
$ws.Range("K61").Value = 1126
$ws.Range("L61").Value = 2050
$ws.Range("M61").Value = -924
$ws.Range("N61").Value = -2454
$ws.Range("H97").Value = 20048
$ws.Range("J97").Value = 20048
$ws.Range("L97").Value = 20048
$ws.Range("N97").Value = -22030
$ws.Range("H113").Value = 1680.4
$ws.Range("I113").Value = 1126
$ws.Range("J113").Value = 2050
$ws.Range("K113").Value = 1126
$ws.Range("L113").Value = 2050
$ws.Range("M113").Value = 1044
$ws.Range("N113").Value = -6390

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").Value = ""
$ws.Range("H62").Value = 2741.8696
$ws.Range("I62").Value = 2704.375
$ws.Range("J62").Value = 2827.5715
$ws.Range("K62").Value = 2704.375
$ws.Range("L62").Value = 2827.5715
$ws.Range("M62").Value = -2080.375
$ws.Range("N62").Value = -4075.5715
$ws.Range("H65").Value = 2741.8696
$ws.Range("I65").Value = 2704.375
$ws.Range("J65").Value = 2827.5715
$ws.Range("K65").Value = 13521.875
$ws.Range("L65").Value = 14137.8575
$ws.Range("M65").Value = -10401.875
$ws.Range("N65").Value = -20377.8575
$ws.Range("H94").Value = 38000
$ws.Range("J94").Value = 38000
$ws.Range("L94").Value = 38000
$ws.Range("N94").Value = -39802
$ws.Range("H97").Value = 48000
$ws.Range("J97").Value = 48000
$ws.Range("L97").Value = 48000
$ws.Range("N97").Value = -49982
$ws.Range("H136").Value = 9215.643
$ws.Range("I136").Value = 21060.8
$ws.Range("K136").Value = 63182.39999999999
$ws.Range("M136").Value = -60632.39999999999
